$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Headers: BTec_Logo-Orange picture, name image1.jpg -> image2.jpg ---
for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $h = $sec.Headers.Item($i)
    if ($h.Exists) {
        $rng = $h.Range
        $xml = $rng.WordOpenXML
        if ($xml -ne $null -and $xml.Contains('name="image1.jpg"')) {
            $newXml = $xml.Replace('name="image1.jpg"', 'name="image2.jpg"')
            $rng.WordOpenXML = $newXml
        }
    }
}

# --- Footers: PearsonLogo picture, name image2.png -> image1.png ---
for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $f = $sec.Footers.Item($i)
    if ($f.Exists) {
        $rng = $f.Range
        $xml = $rng.WordOpenXML
        if ($xml -ne $null -and $xml.Contains('name="image2.png"')) {
            $newXml = $xml.Replace('name="image2.png"', 'name="image1.png"')
            $rng.WordOpenXML = $newXml
        }
    }
}

Write-Output "done"
